# Adds the Artisan Commands `showCurve`, `showExtraCurve`, `showEvents`,
# and `showBackgroundEvents` to the "Commands" sheet of the workbook.
#
# These four new rows are inserted right above the existing "RC Command"
# section (which previously started at row 100), pushing it and everything
# below it down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Drop the trailing empty placeholder row at the very bottom of the sheet
# so it doesn't get dragged along (renumbered) by the row insertion below.
$ws.Rows.Item(1048576).Delete()

# Insert 4 new blank rows above the current row 100 ("RC Command" row),
# matching the row height (13.8) used by the row just above (row 99,
# "keyboard(<bool>)").
$ws.Rows.Item(100).Resize(4).Insert()
$ws.Rows.Item(100).Resize(4).RowHeight = 13.8

$ws.Range("B100").Value = "showCurve(<name>,<bool>)"
$ws.Range("C100").Value = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

$ws.Range("B101").Value = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$ws.Range("C101").Value = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

$ws.Range("B102").Value = "showEvents(<event_type>, <bool>)"
$ws.Range("C102").Value = "shows/hides the events of <event_type> in [1,..,5]"

$ws.Range("B103").Value = "showBackgroundEvents(<bool>)"
$ws.Range("C103").Value = "shows/hides the events of the background profile"

# Match the author's final selection on the Commands sheet.
$ws.Range("C101").Select()
